$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so numeric-looking strings (e.g. "528.21", "1.00")
# are preserved exactly as text rather than being converted to floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.215.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.117.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.21"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.117.26"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.19"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.649.38"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.79"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.48%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "58.233.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.108.53"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.13"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.79"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "344.07"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.516"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.78"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0933"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.40"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.12%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.12"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.65"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.19"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.36"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.25"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.64"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +11.46%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0669"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.156.31"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.83"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.272.24"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.15"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.65"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.19%  "
